$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Picture: add a description ("Alt Text") to the inline screenshot image.
#    This updates the descr="" attribute on both wp:docPr and pic:cNvPr.
# ---------------------------------------------------------------------------
if ($d.InlineShapes.Count -ge 1) {
    $shp = $d.InlineShapes.Item(1)
    $shp.AlternativeText = "Screen-shot of a quotation with a tool-tip showing."
}

# ---------------------------------------------------------------------------
# 2. Wrap the "Put your name and the date..." bullet in a bookmark (as Word
#    does automatically, e.g. when a reviewer copies/pastes or comments on a
#    paragraph). The bookmark covers exactly the paragraph's own text.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Put your name and the date in comments in the head element of the page.*") {
        $r = $p.Range.Duplicate
        [void]$r.MoveEnd(1, 0)
        $d.Bookmarks.Add("_Hlk100049792", $r)
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Wrap the "Test the strings you are passing..." bullet in a bookmark.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Test the strings you are passing*got put together correctly.*") {
        $r = $p.Range.Duplicate
        [void]$r.MoveEnd(1, 0)
        $d.Bookmarks.Add("_Hlk100049812", $r)
        break
    }
}
